$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat is forced to Text ("@") before assignment so that values such as
# "0.160", "6.20", or "0.0000150" keep their original textual formatting instead
# of being auto-coerced into numbers (which would drop trailing/leading zeros).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.685.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.771.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.82%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.33"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.96%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.769.43"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.70%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.78%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -5.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.465"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.46"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.75%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.407.33"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.91%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.774.13"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.628.63"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.31%  "

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.20"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.43%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.26%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "490.29"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.06"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.06%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.743"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.80%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000150"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +12.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.35"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.55%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.37"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.28"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.89%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.60%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.44"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.74"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.09%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.79"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.96%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.329"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "448.55"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.21%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.05"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.33%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.90"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.35"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.20"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -8.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.834.45"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.60%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.92"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0350"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.14"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.21%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.03"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.79%  "
